# edit.ps1 -- apply the "add 2021 spring file" commit:
#   1. Refresh the cached "datetimeFigureOut" footer field from 11/12/20
#      to 11/23/20 on the slide master and on every slide layout.
#   2. On slide 13 ("Class Tortoise: IV"), prefix the first code line of
#      the content placeholder with a purple "const " run so the method
#      signature reads "const int Tortoise::getPosition() const {".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder refresh (slide master + all slide layouts)
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if (-not $sh.HasTextFrame) { continue }
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) { $isDate = $true }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "11/12/20") {
                $tr.Text = "11/23/20"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------
# 2. Slide 13 -- add the "const " keyword run
# ---------------------------------------------------------------------
$slide13 = $p.Slides.Item(13)
$contentShape = $slide13.Shapes.Item(2)
$bodyRange = $contentShape.TextFrame.TextRange
$firstPara = $bodyRange.Paragraphs(1, 1)
$firstRun = $firstPara.Runs(1, 1)

if ($firstRun.Text -eq "int") {
    $constRun = $firstRun.InsertBefore("const ")
    $constRun.Font.Size = $firstRun.Font.Size
    $constRun.Font.Color.RGB = $firstRun.Font.Color.RGB
}
